$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.754521
$ws.Range("H2").Value = 2.263563
$ws.Range("I2").Value = 0.2768403531129761
$ws.Range("J2").Value = 0.2768403531129761
$ws.Range("M2").Value = 11.05182166666667
$ws.Range("N2").Value = 33.155465
$ws.Range("O2").Value = 0.09655996768538078
$ws.Range("P2").Value = 0.09655996768538078
$ws.Range("Q2").Value = 8.338831535755
$ws.Range("R2").Value = 75.04948382179499
$ws.Range("S2").Value = 0.02673169555059838
$ws.Range("T2").Value = 0.02673169555059838

$ws.Range("G3").Value = 0.754521
$ws.Range("H3").Value = 2.263563
$ws.Range("I3").Value = 0.2768403531129761
$ws.Range("J3").Value = 0.2768403531129761
$ws.Range("O3").Value = 0.620602129029037
$ws.Range("P3").Value = 0.620602129029037
$ws.Range("Q3").Value = 53.59463894567499
$ws.Range("R3").Value = 482.351750511075
$ws.Range("S3").Value = 0.1718077125430634
$ws.Range("T3").Value = 0.1718077125430634

$ws.Range("G4").Value = 0.754521
$ws.Range("H4").Value = 2.263563
$ws.Range("I4").Value = 0.2768403531129761
$ws.Range("J4").Value = 0.2768403531129761
$ws.Range("M4").Value = 32.37236033333333
$ws.Range("N4").Value = 97.117081
$ws.Range("O4").Value = 0.2828379032855822
$ws.Range("P4").Value = 0.2828379032855822
$ws.Range("Q4").Value = 24.425625691067
$ws.Range("R4").Value = 219.830631219603
$ws.Range("S4").Value = 0.07830094501931437
$ws.Range("T4").Value = 0.07830094501931437

$ws.Range("G5").Value = 0.9731926666666667
$ws.Range("I5").Value = 0.3570728998754956
$ws.Range("J5").Value = 0.3570728998754956
$ws.Range("M5").Value = 11.05182166666667
$ws.Range("N5").Value = 33.155465
$ws.Range("O5").Value = 0.09655996768538078
$ws.Range("P5").Value = 0.09655996768538078
$ws.Range("Q5").Value = 10.75555179930778
$ws.Range("R5").Value = 96.79996619377
$ws.Range("S5").Value = 0.03447894767330306
$ws.Range("T5").Value = 0.03447894767330306

$ws.Range("G6").Value = 0.9731926666666667
$ws.Range("I6").Value = 0.3570728998754956
$ws.Range("J6").Value = 0.3570728998754956
$ws.Range("O6").Value = 0.620602129029037
$ws.Range("P6").Value = 0.620602129029037
$ws.Range("R6").Value = 622.14462732145
$ws.Range("S6").Value = 0.2216002018813047
$ws.Range("T6").Value = 0.2216002018813047

$ws.Range("G7").Value = 0.9731926666666667
$ws.Range("I7").Value = 0.3570728998754956
$ws.Range("J7").Value = 0.3570728998754956
$ws.Range("M7").Value = 32.37236033333333
$ws.Range("N7").Value = 97.117081
$ws.Range("O7").Value = 0.2828379032855822
$ws.Range("P7").Value = 0.2828379032855822
$ws.Range("Q7").Value = 31.50454367909089
$ws.Range("R7").Value = 283.540893111818
$ws.Range("S7").Value = 0.1009937503208878
$ws.Range("T7").Value = 0.1009937503208878

$ws.Range("G8").Value = 0.7824410000000001
$ws.Range("H8").Value = 2.347323
$ws.Range("I8").Value = 0.2870844452706686
$ws.Range("J8").Value = 0.2870844452706686
$ws.Range("M8").Value = 11.05182166666667
$ws.Range("N8").Value = 33.155465
$ws.Range("O8").Value = 0.09655996768538078
$ws.Range("P8").Value = 0.09655996768538078
$ws.Range("Q8").Value = 8.647398396688335
$ws.Range("R8").Value = 77.82658557019501
$ws.Range("S8").Value = 0.02772086475831122
$ws.Range("T8").Value = 0.02772086475831123

$ws.Range("G9").Value = 0.7824410000000001
$ws.Range("H9").Value = 2.347323
$ws.Range("I9").Value = 0.2870844452706686
$ws.Range("J9").Value = 0.2870844452706686
$ws.Range("O9").Value = 0.620602129029037
$ws.Range("P9").Value = 0.620602129029037
$ws.Range("Q9").Value = 55.57783400500833
$ws.Range("R9").Value = 500.200506045075
$ws.Range("S9").Value = 0.178165217946097
$ws.Range("T9").Value = 0.178165217946097

$ws.Range("G10").Value = 0.7824410000000001
$ws.Range("H10").Value = 2.347323
$ws.Range("I10").Value = 0.2870844452706686
$ws.Range("J10").Value = 0.2870844452706686
$ws.Range("M10").Value = 32.37236033333333
$ws.Range("N10").Value = 97.117081
$ws.Range("O10").Value = 0.2828379032855822
$ws.Range("P10").Value = 0.2828379032855822
$ws.Range("Q10").Value = 25.32946199157367
$ws.Range("R10").Value = 227.965157924163
$ws.Range("S10").Value = 0.08119836256626038
$ws.Range("T10").Value = 0.0811983625662604

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.2153186666666667
$ws.Range("H11").Value = 0.645956
$ws.Range("I11").Value = 0.07900230174085969
$ws.Range("J11").Value = 0.07900230174085969
$ws.Range("M11").Value = 11.05182166666667
$ws.Range("N11").Value = 33.155465
$ws.Range("O11").Value = 0.09655996768538078
$ws.Range("P11").Value = 0.09655996768538078
$ws.Range("Q11").Value = 2.379663505504444
$ws.Range("R11").Value = 21.41697154954
$ws.Range("S11").Value = 0.007628459703168113
$ws.Range("T11").Value = 0.007628459703168113

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.2153186666666667
$ws.Range("H12").Value = 0.645956
$ws.Range("I12").Value = 0.07900230174085969
$ws.Range("J12").Value = 0.07900230174085969
$ws.Range("O12").Value = 0.620602129029037
$ws.Range("P12").Value = 0.620602129029037
$ws.Range("Q12").Value = 15.29437377921111
$ws.Range("R12").Value = 137.6493640129
$ws.Range("S12").Value = 0.04902899665857192
$ws.Range("T12").Value = 0.04902899665857192

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.2153186666666667
$ws.Range("H13").Value = 0.645956
$ws.Range("I13").Value = 0.07900230174085969
$ws.Range("J13").Value = 0.07900230174085969
$ws.Range("M13").Value = 32.37236033333333
$ws.Range("N13").Value = 97.117081
$ws.Range("O13").Value = 0.2828379032855822
$ws.Range("P13").Value = 0.2828379032855822
$ws.Range("Q13").Value = 6.970373463826222
$ws.Range("R13").Value = 62.733361174436
$ws.Range("S13").Value = 0.02234484537911966
$ws.Range("T13").Value = 0.02234484537911966
